# Remove two unnecessary bypass capacitors
#
# Row 4 (1 uF caps): drop C31 and C35 -> qty 13 -> 11
# Row 3 (0.1 uF caps): drop C30 and C33 -> qty 14 -> 12
#
# NOTE: B4 is updated before B3 so that the regenerated shared-strings
# table appends the new "C6, C8, ..." text ahead of the new
# "C2, C3, ..." text (matches upstream ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "C6, C8, C9, C10, C13, C14, C20, C23, C26, C28, C34"
$ws.Range("H4").Value = 11

$ws.Range("B3").Value = "C2, C3, C4, C5, C12, C15, C18, C21, C24, C25, C29, C32"
$ws.Range("H3").Value = 12

# Column B got a bit narrower after the reference lists shrank.
$ws.Range("B:B").ColumnWidth = 44.15

# Selection/scroll position left where the author was working.
$ws.Range("H5").Select()
